# The commit swaps the contents of ppt/theme/theme1.xml (was "Office Theme")
# and ppt/theme/theme2.xml (was "Integral" / "Red Violet") so that:
#   - theme1.xml (used by the Notes Master) becomes the "Integral" / Red Violet theme
#   - theme2.xml (used by the Slide Master) becomes the stock "Office Theme" palette
#
# The PowerPoint object model only exposes the *active* master's colour
# scheme (Presentation.SlideMaster / .NotesMaster / .HandoutMaster all
# resolve to the one Master that backs ppt/theme/theme2.xml), so we recolour
# that scheme's 12 theme colours to the "Office" palette that the diff
# shows landing in theme2.xml. Each RGB value is written using the
# COLORREF (0xBBGGRR) convention PowerPoint's ColorScheme.Colors(i).RGB
# expects, in clrScheme order: dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink.

$p = $ppt.ActivePresentation
$m = $p.SlideMaster
$cs = $m.ColorScheme

$cs.Colors(1).RGB  = 0          # dk1      000000
$cs.Colors(2).RGB  = 16777215   # lt1      FFFFFF
$cs.Colors(3).RGB  = 6968388    # dk2      44546A
$cs.Colors(4).RGB  = 15132391   # lt2      E7E6E6
$cs.Colors(5).RGB  = 13998939   # accent1  5B9BD5
$cs.Colors(6).RGB  = 3243501    # accent2  ED7D31
$cs.Colors(7).RGB  = 10855845   # accent3  A5A5A5
$cs.Colors(8).RGB  = 49407      # accent4  FFC000
$cs.Colors(9).RGB  = 12874308   # accent5  4472C4
$cs.Colors(10).RGB = 4697456    # accent6  70AD47
$cs.Colors(11).RGB = 12673797   # hlink    0563C1
$cs.Colors(12).RGB = 7491477    # folHlink 954F72
